$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fitnessValues = @(13738,11350,11350,9970,9970,9970,9970,9970,9913,9913,9878,9878,9878,9308,9308,9308,9308,9241,9241,9241,9241,9166,9166,9159,9159,8951,8951,8951,8375,8375,8375,8375,8375,8375,8375,8375,8375,8375,8375,8375,8375,8375,8375,8375,8342,8283,8283,8283,8283,8283,8283,8225,8225,8225,8225,8225,8225,8225,8225,8197,8197,8197,8197,8197,8197,8197,8070,8070,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7651,7619,7619,7619,7619,7619,7619,7619,7619,7619,7619,7619,7619,7619,7619,7619,7619,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569)

$startRow = 2
for ($i = 0; $i -lt $fitnessValues.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $fitnessValues[$i]
}
